# The "Debugging" slide and the "Live Demo" slide that follows it were saved
# in the wrong order. Fix the slide order by moving the "Debugging" slide
# (currently slide 15) so it comes right after the "Live Demo" slide
# (currently slide 16) - i.e. swap slides 15 and 16.

$p = $ppt.ActivePresentation

$debuggingIndex = -1
$liveDemoIndex = -1

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
    if ($title -eq "Debugging" -and $debuggingIndex -eq -1) {
        $debuggingIndex = $i
    }
    if ($title -eq "Live Demo" -and $debuggingIndex -ne -1 -and $liveDemoIndex -eq -1 -and $i -eq ($debuggingIndex + 1)) {
        $liveDemoIndex = $i
    }
}

if ($debuggingIndex -eq -1 -or $liveDemoIndex -eq -1) {
    # Fall back to the known positions (15 and 16) if the titles could not
    # be matched for some reason.
    $debuggingIndex = 15
    $liveDemoIndex = 16
}

# Move the "Debugging" slide to just after the "Live Demo" slide, which
# swaps the two adjacent slides so "Live Demo" now comes first.
$p.Slides.Item($debuggingIndex).MoveTo($liveDemoIndex)
